$d = $word.ActiveDocument

# --- First line: reuse the existing (only) paragraph ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.LanguageID = "ka-GE"
$r1.InsertAfter("პირველი სტრიქონი")
$r1.Font.Name = "Sylfaen"
$r1.LanguageID = "ka-GE"

# --- Second line: new paragraph appended after the first ---
$r1.InsertParagraphAfter()

$p2 = $d.Paragraphs(2)
$r2 = $p2.Range
$r2.LanguageID = "ka-GE"
$r2.InsertAfter("მეორე სტრიქონი")
$r2.Font.Name = "Sylfaen"
$r2.LanguageID = "ka-GE"

Write-Output "done"
